$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextCell $ws 'D2' '65.954.45'
$ws.Range('E2').Value = '  -0.90%  '
Set-TextCell $ws 'D3' '3.287.86'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextCell $ws 'D5' '584.82'
$ws.Range('E5').Value = '  +2.34%  '
Set-TextCell $ws 'D6' '180.52'
$ws.Range('E6').Value = '  -0.55%  '
Set-TextCell $ws 'D7' '0.641'
$ws.Range('E7').Value = '  +7.17%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -3.17%  '
Set-TextCell $ws 'D10' '6.77'
$ws.Range('E10').Value = '  +2.38%  '
Set-TextCell $ws 'D11' '0.403'
$ws.Range('E11').Value = '  +0.10%  '
Set-TextCell $ws 'D12' '3.853.26'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('E13').Value = '  -4.52%  '
Set-TextCell $ws 'D14' '65.980.81'
$ws.Range('E14').Value = '  -0.97%  '
Set-TextCell $ws 'D15' '26.15'
$ws.Range('E15').Value = '  -3.54%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D16' '3.354.31'
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D17' '0.0000163'
$ws.Range('E17').Value = '  -2.71%  '
Set-TextCell $ws 'D18' '426.59'
$ws.Range('E18').Value = '  -1.00%  '
Set-TextCell $ws 'D19' '13.22'
$ws.Range('E19').Value = '  -3.05%  '
Set-TextCell $ws 'D20' '5.50'
$ws.Range('E20').Value = '  -3.24%  '
Set-TextCell $ws 'D21' '7.37'
$ws.Range('E21').Value = '  -3.04%  '
Set-TextCell $ws 'D22' '71.74'
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('E23').Value = '  +0.12%  '
Set-TextCell $ws 'D24' '5.69'
$ws.Range('E24').Value = '  +0.40%  '
Set-TextCell $ws 'D25' '3.420.28'
$ws.Range('E25').Value = '  -0.87%  '
Set-TextCell $ws 'D26' '0.511'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('E28').Value = '  -4.52%  '
Set-TextCell $ws 'D29' '8.95'
$ws.Range('E29').Value = '  -0.98%  '
Set-TextCell $ws 'D30' '0.999'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -0.37%  '
Set-TextCell $ws 'D32' '22.25'
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('E33').Value = '  +0.06%  '
Set-TextCell $ws 'D34' '5.17'
$ws.Range('E34').Value = '  -2.42%  '
Set-TextCell $ws 'D35' '6.59'
$ws.Range('E35').Value = '  -2.70%  '
Set-TextCell $ws 'D36' '1.19'
$ws.Range('E36').Value = '  -3.54%  '
Set-TextCell $ws 'D37' '158.96'
$ws.Range('E37').Value = '  -0.62%  '
Set-TextCell $ws 'D38' '1.44'
$ws.Range('E38').Value = '  -3.76%  '
Set-TextCell $ws 'D39' '1.81'
$ws.Range('E39').Value = '  -1.68%  '
Set-TextCell $ws 'D40' '26.58'
$ws.Range('E40').Value = '  -2.04%  '
Set-TextCell $ws 'D41' '2.818.91'
$ws.Range('E41').Value = '  +0.17%  '
Set-TextCell $ws 'D42' '0.764'
Set-TextCell $ws 'D43' '4.32'
$ws.Range('E43').Value = '  -2.52%  '
Set-TextCell $ws 'D44' '40.02'
$ws.Range('E44').Value = '  -0.28%  '
Set-TextCell $ws 'D45' '0.0660'
$ws.Range('E45').Value = '  -2.12%  '
Set-TextCell $ws 'D46' '5.94'
$ws.Range('E46').Value = '  -4.06%  '
Set-TextCell $ws 'D47' '2.29'
$ws.Range('E47').Value = '  -2.36%  '
Set-TextCell $ws 'D48' '314.66'
$ws.Range('E48').Value = '  -1.76%  '
Set-TextCell $ws 'D49' '23.15'
$ws.Range('E49').Value = '  -4.94%  '
Set-TextCell $ws 'D50' '0.0268'
$ws.Range('E50').Value = '  -1.27%  '
$ws.Range('E51').Value = '  +3.57%  '

Write-Host "Applied all cell updates."